# Fruta / hortaliza, semanal
# Insert two new data rows (new weekly observation) above the existing
# row 52, shifting the rest of the table (old rows 52-176) down by two
# rows to 54-178, and populate the two new rows with the latest figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 52:176 down to 54:178, carrying formatting along.
$ws.Rows("52:53").Insert()

# New row 52: Primera quality observation for the new date.
$ws.Range("A52").Value = 11
$ws.Range("B52").Value = "Vega Monumental Concepción"
$ws.Range("C52").Value = "Bíobío"
$ws.Range("D52").Value = 44581
$ws.Range("E52").Value = 8
$ws.Range("F52").Value = 100114013
$ws.Range("G52").Value = "Zanahoria"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 700
$ws.Range("K52").Value = 8000
$ws.Range("L52").Value = 8500
$ws.Range("M52").Value = 8286
$ws.Range("N52").Value = "$/saco 20 kilos"
$ws.Range("O52").Value = "Región de Ñuble"
$ws.Range("P52").Value = 414
$ws.Range("Q52").Value = 20
$ws.Range("R52").Value = "Hortaliza"

# New row 53: Segunda quality observation for the new date.
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = "Vega Monumental Concepción"
$ws.Range("C53").Value = "Bíobío"
$ws.Range("D53").Value = 44581
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = 100114013
$ws.Range("G53").Value = "Zanahoria"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Segunda"
$ws.Range("J53").Value = 400
$ws.Range("K53").Value = 7000
$ws.Range("L53").Value = 7000
$ws.Range("M53").Value = 7000
$ws.Range("N53").Value = "$/saco 20 kilos"
$ws.Range("O53").Value = "Región de Ñuble"
$ws.Range("P53").Value = 350
$ws.Range("Q53").Value = 20
$ws.Range("R53").Value = "Hortaliza"
